{"js": "// Lattice multiplication exercises: replace the 15 practice problems\n// (5 rows x 3 columns) in the first table with a new set of problems,\n// keeping the existing \"lattice grid\" template/formatting intact.\n\n// New problems, in row-major order (5 rows x 3 cols), as [multiplicand, multiplier]\n// pairs - both two-digit numbers, matching the cells being replaced.\nconst problems = [\n  [\"55\", \"87\"], [\"24\", \"95\"], [\"55\", \"19\"],\n  [\"35\", \"35\"], [\"31\", \"97\"], [\"51\", \"80\"],\n  [\"39\", \"31\"], [\"72\", \"51\"], [\"24\", \"26\"],\n  [\"67\", \"12\"], [\"41\", \"65\"], [\"50\", \"17\"],\n  [\"80\", \"24\"], [\"99\", \"57\"], [\"46\", \"19\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = 3;\n\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const idx = r * colCount + c;\n    if (idx >= problems.length) continue;\n    const [a, b] = problems[idx];\n\n    // Lines of the lattice-multiplication \"card\":\n    //   \"A x B\"\n    //   \"  B0    B1\"      (digits of the multiplier, spaced out)\n    //   \"  ----\"\n    //   \"A0|    |\"        (tens digit of the multiplicand)\n    //   \"A1|    |\"        (ones digit of the multiplicand)\n    const lines = [\n      `${a} x ${b}`,\n      `  ${b[0]}    ${b[1]}`,\n      \"  ----\",\n      `${a[0]}|    |`,\n      `${a[1]}|    |`,\n    ];\n    // \\v (vertical tab) becomes a line break (<w:br/>) when inserted via Office.js.\n    const newText = lines.join(\"\\v\");\n\n    const cell = table.getCell(r, c);\n    cell.body.paragraphs.load(\"items\");\n    await context.sync();\n\n    const para = cell.body.paragraphs.items[0];\n    const rng = para.getRange();\n    // Replace the whole paragraph range in one shot so the existing run\n    // formatting (sz=32) on the first run is reused instead of rebuilt.\n    rng.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Lattice multiplication exercises: replace the 15 practice problems\n# (5 rows x 3 columns) in the first table with a new set of problems,\n# keeping the existing \"lattice grid\" template/formatting intact.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# New problems, in row-major order (5 rows x 3 cols), as two-digit\n# [multiplicand, multiplier] string pairs.\n$problems = @(\n    @(\"55\", \"87\"), @(\"24\", \"95\"), @(\"55\", \"19\"),\n    @(\"35\", \"35\"), @(\"31\", \"97\"), @(\"51\", \"80\"),\n    @(\"39\", \"31\"), @(\"72\", \"51\"), @(\"24\", \"26\"),\n    @(\"67\", \"12\"), @(\"41\", \"65\"), @(\"50\", \"17\"),\n    @(\"80\", \"24\"), @(\"99\", \"57\"), @(\"46\", \"19\")\n)\n\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n$vt = [char]11\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($idx -ge $problems.Count) { continue }\n        $pair = $problems[$idx]\n        $a = $pair[0]\n        $b = $pair[1]\n\n        # Lines of the lattice-multiplication \"card\":\n        #   \"A x B\"\n        #   \"  B0    B1\"      (digits of the multiplier, spaced out)\n        #   \"  ----\"\n        #   \"A0|    |\"        (tens digit of the multiplicand)\n        #   \"A1|    |\"        (ones digit of the multiplicand)\n        $line1 = \"{0} x {1}\" -f $a, $b\n        $line2 = \"  {0}    {1}\" -f $b[0], $b[1]\n        $line3 = \"  ----\"\n        $line4 = \"{0}|    |\" -f $a[0]\n        $line5 = \"{0}|    |\" -f $a[1]\n\n        $text = $line1 + $vt + $line2 + $vt + $line3 + $vt + $line4 + $vt + $line5\n\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $text\n\n        $idx++\n    }\n}\n"}
